$d = $word.ActiveDocument
$targetRange = $d.Content
$targetRange.Collapse(0)

$xmlPayload = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>Includ</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>e command</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>The instruction “</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>include</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">”  </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>allows</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> us to include another file inside our </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>PHP</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> file. </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>One of th</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>e</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> instruction</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>’</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">s applications would be having </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>the</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> same header and footer.</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> In this way, we can write a separate file for </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">header and another one for </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">footer, then we can just simply include them </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>n every other page.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $targetRange.InsertXML($xmlPayload)

Write-Output "Inserted Include command section."
